# Sprint 1 Backlog Burndown - update sprint 1 backlog data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: "Login/Logout" task was reassigned to Jacob Evans, and its
# remaining-work tracking was updated for Week 1/Week 2.
$ws.Range("C6").Value = "Jacob Evans"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 4

# Row 11 (Landing Page task): fill in Week 1 / Week 2 remaining amounts.
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0

# Row 12 (User Management Page task): fill in Week 1 / Week 2 remaining amounts.
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 0

# Row 15 (Group Management Page task): fill in Week 1 / Week 2 remaining amounts.
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 4

# Row 22 (Database task): fill in Week 1 / Week 2 remaining amounts.
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 0

# Row 23 (entity relationships task): fill in Week 1 / Week 2 remaining amounts.
$ws.Range("E23").Value = 6
$ws.Range("F23").Value = 0

# Recalculate so the Estimate Totals row (and the burndown chart that
# references it) reflect the updated weekly totals.
$excel.CalculateFullRebuild()

# Update the active selection / scroll position to match where the user left
# off reviewing the backlog.
$null = $ws.Range("E12").Select()
